$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9549548625946045
$ws.Range("B1").Value = 1.4534672498703
$ws.Range("C1").Value = 2.715541124343872
$ws.Range("D1").Value = 1.568399310112
$ws.Range("E1").Value = 1.446734309196472
